$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells C4 and C11 hold the "date" as literal text (e.g. "06.03.2023").
# Assigning a plain string via .Value would make Excel auto-detect it as a
# real date and reformat/retype the cell. To keep the cell as plain text
# (same style, same inline/shared-string type) we build the text via a
# formula and then paste back just the computed value.
$c4 = $ws.Range("C4")
$c4.Formula = '="07.03.2023"'
$c4.Copy()
$c4.PasteSpecial(-4163) | Out-Null

$c11 = $ws.Range("C11")
$c11.Formula = '="07.03.2023"'
$c11.Copy()
$c11.PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0
